$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already has two sample rows (2 and 3) with identical
# formatting/styles. Add a third data row ("TI" / 144587410 / ...),
# reusing the row-3 formatting so every cell keeps the same style index.
$ws.Range("A3:P3").Copy() | Out-Null
$ws.Range("A4:P4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A4").Value = "TI"
$ws.Range("B4").Value = 144587410
$ws.Range("C4").Value = 44552
$ws.Range("D4").Value = 2021
$ws.Range("E4").Value = "DIC"
$ws.Range("F4").Value = "Consulta Retina"
$ws.Range("G4").Value = "SEBASTIAN ROJAS MUNERA"
$ws.Range("H4").Value = "H353"
$ws.Range("I4").Value = "DEGENERACION DE LA MACULA Y DEL POLO POSTERIOR DEL OJO"
$ws.Range("J4").Value = "OI"
$ws.Range("K4").Value = "20/150"
$ws.Range("L4").Value = "CC"
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = "SI"
$ws.Range("O4").Value = ""
$ws.Range("P4").Value = "EPS Y MEDICINA PREPAGADA SURAMERICANA S.A"

# Row height of the new row settles a bit smaller than the others once
# it is filled in.
$ws.Rows.Item(4).RowHeight = 13.8

# The "no duplicates" conditional format over column B now has to cover
# the new row too.
$fc = $ws.Range("B2:B3").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("B2:B4"))

# Column C (the date column) is widened to fit the new content.
$ws.Columns.Item(3).ColumnWidth = 16.33

# Leave the selection where the user ended up after entering the row.
$ws.Range("D6").Select() | Out-Null
